$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the "Updated symbol list" commit diff.
# Column D holds numeric-looking price strings stored as text (inlineStr) in the
# original workbook; a leading apostrophe forces Excel to keep them as text
# instead of silently converting them to a Number cell (matches typing them in
# the Excel UI with a text/quote prefix).

# Row 2
$ws.Range('D2').Value = "'247.74"

# Row 4
$ws.Range('D4').Value = "'5.550"

# Row 5
$ws.Range('D5').Value = "'0.05639"

# Row 7
$ws.Range('D7').Value = "'1.072"

# Row 8
$ws.Range('D8').Value = "'0.8017"

# Row 9
$ws.Range('B9').Value = 'One'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D9').Value = "'0.01172"
$ws.Range('E9').Value = '8OneONEBestin24h'

# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'0.1427"
$ws.Range('E10').Value = '9WazirXWRX'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = "'0.07318"
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = "'0.03198"
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.02993"
$ws.Range('E13').Value = '12BitrueCoinBTR'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09267"
$ws.Range('E14').Value = '13BitMartTokenBMX'

# Row 15
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001662"
$ws.Range('E15').Value = '14BitForexTokenBF'

# Row 16
$ws.Range('D16').Value = "'0.04690"

# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = "'0.006266"
$ws.Range('E17').Value = '16TigerCashTCH'

# Row 18
$ws.Range('B18').Value = 'BitKan'
$ws.Range('C18').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D18').Value = "'0.001054"
$ws.Range('E18').Value = '17BitKanKAN'

# Row 19
$ws.Range('B19').Value = 'HotbitToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D19').Value = "'0.003833"
$ws.Range('E19').Value = '18HotbitTokenHTB'

# Row 20
$ws.Range('B20').Value = 'NitroEx'
$ws.Range('C20').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D20').Value = "'0.0001501"
$ws.Range('E20').Value = '19NitroExNTX'

# Row 21
$ws.Range('B21').Value = 'UpBots'
$ws.Range('C21').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D21').Value = "'0.0004003"
$ws.Range('E21').Value = '20UpBotsUBXT'

# Row 22
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').Value = "'3.982"
$ws.Range('E22').Value = '21LEOLEO'

# Row 23
$ws.Range('B23').Value = 'GateToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D23').Value = "'3.397"
$ws.Range('E23').Value = '22GateTokenGT'

# Row 24
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D24').Value = "'2.098"
$ws.Range('E24').Value = '23BTSETokenBTSE'

# Row 25
$ws.Range('B25').Value = 'BitpandaEcosystemToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D25').Value = "'0.3290"
$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'

# Row 26
$ws.Range('B26').Value = 'ProBitToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D26').Value = "'0.1292"
$ws.Range('E26').Value = '25ProBitTokenPROB'

# Row 27
$ws.Range('B27').Value = 'MCDex'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D27').Value = "'2.585"
$ws.Range('E27').Value = '26MCDexMCB'

# Row 40
$ws.Range('D40').Value = "'0.04196"

# Row 41
$ws.Range('D41').Value = "'0.007026"

# Row 42
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = "'0.003503"
$ws.Range('E42').Value = '41CEJICEJI'

# Row 43
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').Value = "'0.1046"
$ws.Range('E43').Value = '42BKEXTokenBKK'

# Row 44
$ws.Range('D44').Value = "'0.008734"

# Row 45
$ws.Range('D45').Value = "'0.00005635"

# Row 48
$ws.Range('D48').Value = "'0.02720"
